# Fruta / hortaliza, semanal
# Weekly refresh of the Membrillo (quince) price series for Vega Monumental
# Concepcion: existing observations are re-dated/re-priced and two more
# weekly records (rows 15-16) are appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-12) with revised values ---
# Row 2
$ws.Range("D2").Value = 44698
$ws.Range("Q2").Value = '$/caja 18 kilos granel'

# Row 3
$ws.Range("D3").Value = 44776
$ws.Range("L3").Value = 'Primera'
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("S3").Value = 556

# Row 4
$ws.Range("D4").Value = 44776
$ws.Range("L4").Value = 'Segunda'
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("S4").Value = 444

# Row 5
$ws.Range("D5").Value = 44999
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("S5").Value = 667

# Row 6
$ws.Range("D6").Value = 44999
$ws.Range("L6").Value = 'Segunda'
$ws.Range("N6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("S6").Value = 556
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44307
$ws.Range("M7").Value = 50
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 556

# Row 8
$ws.Range("D8").Value = 44307
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("Q8").Value = '$/bandeja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 444

# Row 9
$ws.Range("D9").Value = 44316
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 9500
$ws.Range("S9").Value = 528

# Row 10
$ws.Range("D10").Value = 44358
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 11500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 639

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("M11").Value = 100
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 10500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Región del Maule'
$ws.Range("S11").Value = 583

# Row 12
$ws.Range("D12").Value = 44299
$ws.Range("N12").Value = 9000
$ws.Range("O12").Value = 9000
$ws.Range("P12").Value = 9000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región del Maule'
$ws.Range("S12").Value = 500

# --- Append new rows 15 and 16 (weekly data) ---
# Row 15
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 'Vega Monumental Concepción'
$ws.Range("C15").Value = 'Bíobío'
$ws.Range("D15").Value = 44363
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100104
$ws.Range("H15").Value = 'Frutos de pepita'
$ws.Range("I15").Value = 100104003
$ws.Range("J15").Value = 'Membrillo'
$ws.Range("K15").Value = 'Champion'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 9500
$ws.Range("Q15").Value = '$/caja 15 kilos empedrada'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 633
$ws.Range("T15").Value = 15
$ws.Range("D15").NumberFormat = $ws.Range("D2").NumberFormat

# Row 16
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 'Vega Monumental Concepción'
$ws.Range("C16").Value = 'Bíobío'
$ws.Range("D16").Value = 44425
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = 'Frutos de pepita'
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = 'Membrillo'
$ws.Range("K16").Value = 'Champion'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 13000
$ws.Range("P16").Value = 12500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 694
$ws.Range("T16").Value = 18
$ws.Range("D16").NumberFormat = $ws.Range("D2").NumberFormat
